# Insert two new data rows for Albahaca (Vega Central Mapocho de Santiago)
# right before the current row 547, shifting existing rows 547-575 down to 549-577.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 547 (pushes old 547.. down by 2)
$ws.Rows.Item(547).Insert()
$ws.Rows.Item(547).Insert()

# New row 547: Primera quality, 2023-04-25 (serial 45041)
$ws.Cells.Item(547, 1).Value = 9
$ws.Cells.Item(547, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(547, 3).Value = "Metropolitana"
$ws.Cells.Item(547, 4).Value = 45041
$ws.Cells.Item(547, 5).Value = 13
$ws.Cells.Item(547, 6).Value = 100112052
$ws.Cells.Item(547, 7).Value = "Albahaca"
$ws.Cells.Item(547, 8).Value = "Sin especificar"
$ws.Cells.Item(547, 9).Value = "Primera"
$ws.Cells.Item(547, 10).Value = 340
$ws.Cells.Item(547, 11).Value = 3000
$ws.Cells.Item(547, 12).Value = 3000
$ws.Cells.Item(547, 13).Value = 3000
$ws.Cells.Item(547, 14).Value = "$/docena de matas"
$ws.Cells.Item(547, 15).Value = "Región Metropolitana"
$ws.Cells.Item(547, 16).Value = 500
$ws.Cells.Item(547, 17).Value = 6
$ws.Cells.Item(547, 18).Value = "Hortaliza"

# New row 548: Segunda quality, same date 2023-04-25 (serial 45041)
$ws.Cells.Item(548, 1).Value = 9
$ws.Cells.Item(548, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(548, 3).Value = "Metropolitana"
$ws.Cells.Item(548, 4).Value = 45041
$ws.Cells.Item(548, 5).Value = 13
$ws.Cells.Item(548, 6).Value = 100112052
$ws.Cells.Item(548, 7).Value = "Albahaca"
$ws.Cells.Item(548, 8).Value = "Sin especificar"
$ws.Cells.Item(548, 9).Value = "Segunda"
$ws.Cells.Item(548, 10).Value = 160
$ws.Cells.Item(548, 11).Value = 2500
$ws.Cells.Item(548, 12).Value = 2500
$ws.Cells.Item(548, 13).Value = 2500
$ws.Cells.Item(548, 14).Value = "$/docena de matas"
$ws.Cells.Item(548, 15).Value = "Región Metropolitana"
$ws.Cells.Item(548, 16).Value = 417
$ws.Cells.Item(548, 17).Value = 6
$ws.Cells.Item(548, 18).Value = "Hortaliza"
